$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.93'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.17'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.341'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05734'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.138'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8165'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8670'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1371'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07001'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03219'
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1331'
$ws.Range("E13").Value = '12ProBitTokenPROB'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.02882'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09374'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.741'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001524'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04691'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0005989'
$ws.Range("E19").Value = '18OneONE'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006160'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001239'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.004785'
$ws.Range("E22").Value = '21HotbitTokenHTB'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.00006799'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.529'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.151'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3170'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03702'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006449'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002220'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008626'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005478'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002676'
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
